$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old contents entirely; we'll rewrite the whole table with the
# new room/feature layout (the Elevator room was split into three rows and a
# new "Features (must have two)" column of data was added under column D).
$ws.Cells.Clear()

$data = [ordered]@{
    1  = [ordered]@{ A = "Room";                        B = "Person"; C = "Items";               D = "Features (must have two)" }
    2  = [ordered]@{ A = "Bridge";                       B = "Robert"; C = "button";              D = "window" }
    3  = [ordered]@{                                                                              D = "notepad" }
    4  = [ordered]@{ A = "Elevator 1";                   B = "Robert";                            D = "map" }
    5  = [ordered]@{                                                                              D = "speaker" }
    6  = [ordered]@{ A = "Elevator 2";                   B = "Robert";                            D = "map" }
    7  = [ordered]@{                                                                              D = "speaker" }
    8  = [ordered]@{ A = "Elevator 3";                   B = "Robert";                            D = "map" }
    9  = [ordered]@{                                                                              D = "speaker" }
    10 = [ordered]@{ A = "Crew Quarters";                B = "Robert";                            D = "diary" }
    11 = [ordered]@{                                                                              D = "pajamas" }
    12 = [ordered]@{ A = "Shuttle Bay";                  B = "Robert" }
    14 = [ordered]@{ A = "Captain's Room";               B = "Robert"; C = "model sailing ship" }
    16 = [ordered]@{ A = "Transporter Room";             B = "James" }
    18 = [ordered]@{ A = "Engineering Access Tube";      B = "James" }
    20 = [ordered]@{ A = "Auxiliary Bridge";             B = "James" }
    22 = [ordered]@{ A = "Main Engineering";             B = "James" }
    24 = [ordered]@{ A = "Engineering Core";             B = "James"; C = "reactor fuel" }
    26 = [ordered]@{ A = "Armory";                       B = "James"; C = "blaster pistol";       D = "pistol instructions" }
    27 = [ordered]@{                                                                              D = "bow" }
    28 = [ordered]@{ A = "Forward Observation Lounge";   B = "Brent";                             D = "table of notes" }
    29 = [ordered]@{                                                                              D = "telescope" }
    30 = [ordered]@{ A = "Medical Bay";                  B = "Brent"; C = "android";              D = "dying alien" }
    31 = [ordered]@{                                                                              D = "alien notes" }
    32 = [ordered]@{ A = "Mess Hall";                    B = "Brent"; C = "plastic pass key";     D = "leftovers" }
    33 = [ordered]@{                                                                              D = "dying man" }
    34 = [ordered]@{ A = "Arboretum";                    B = "Brent"; C = "blue rose";            D = "crops" }
    35 = [ordered]@{                                                                              D = "turtle" }
    36 = [ordered]@{ A = "Cargo Deck";                   B = "Brent"; C = "hibernation pod";      D = "manifest" }
    37 = [ordered]@{                                                                              D = "packages" }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

# Bold header row style (same as before the edit).
$ws.Range("A1:D1").Style = "Bold"

# Restore the selected cell recorded in the saved workbook.
$ws.Range("D11").Select()
